$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 496.25
$ws.Range("I33").Value = 263.85715
$ws.Range("K33").Value = 263.85715
$ws.Range("M33").Value = -34.85714999999999
$ws.Range("H41").Value = 568.75
$ws.Range("I41").Value = 78
$ws.Range("J41").Value = 2041
$ws.Range("K41").Value = 78
$ws.Range("L41").Value = 2041
$ws.Range("M41").Value = 362
$ws.Range("N41").Value = -2921
$ws.Range("H62").Value = 20812.25
$ws.Range("I62").Value = 18600.2
$ws.Range("K62").Value = 18600.2
$ws.Range("M62").Value = -17976.2
$ws.Range("H64").Value = 7423.722
$ws.Range("I64").Value = 4089.4
$ws.Range("K64").Value = 4089.4
$ws.Range("M64").Value = -3841.4
$ws.Range("H65").Value = 20812.25
$ws.Range("I65").Value = 18600.2
$ws.Range("K65").Value = 93001
$ws.Range("M65").Value = -89881
$ws.Range("H67").Value = 7423.722
$ws.Range("I67").Value = 4089.4
$ws.Range("K67").Value = 4089.4
$ws.Range("M67").Value = -3231.4
$ws.Range("H70").Value = 2997.7144
$ws.Range("I70").Value = 2300
$ws.Range("J70").Value = 3276.8
$ws.Range("K70").Value = 6900
$ws.Range("L70").Value = 9830.400000000001
$ws.Range("M70").Value = -6630
$ws.Range("N70").Value = -10370.4
$ws.Range("H73").Value = 2997.7144
$ws.Range("I73").Value = 2300
$ws.Range("J73").Value = 3276.8
$ws.Range("K73").Value = 6900
$ws.Range("L73").Value = 9830.400000000001
$ws.Range("M73").Value = -5964
$ws.Range("N73").Value = -11702.4
$ws.Range("H113").Value = 4314.2
$ws.Range("I113").Value = 3305
$ws.Range("J113").Value = 4566.5
$ws.Range("K113").Value = 3305
$ws.Range("L113").Value = 4566.5
$ws.Range("M113").Value = -51
$ws.Range("N113").Value = -11074.5
$ws.Range("H138").Value = 1760.9412
$ws.Range("I138").Value = 1120.2916
$ws.Range("J138").Value = 3298.5
$ws.Range("K138").Value = 3360.8748
$ws.Range("L138").Value = 9895.5
$ws.Range("M138").Value = 1779.1252
$ws.Range("N138").Value = -20175.5
$ws.Range("H141").Value = 1839.6666
$ws.Range("I141").Value = 1643.2727
$ws.Range("K141").Value = 4929.8181
$ws.Range("M141").Value = 250.1818999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5955321.5
$ws.Range("I32").Value = 5955321.5
$ws.Range("K32").Value = 5955321.5
$ws.Range("M32").Value = -5955034.5
$ws.Range("H45").Value = 3199.7144
$ws.Range("I45").Value = 3039.6
$ws.Range("J45").Value = 3600
$ws.Range("K45").Value = 3039.6
$ws.Range("L45").Value = 3600
$ws.Range("M45").Value = -2662.6
$ws.Range("N45").Value = -4354
$ws.Range("H61").Value = 1755846.4
$ws.Range("I61").Value = 2382156.5
$ws.Range("K61").Value = 2382156.5
$ws.Range("M61").Value = -2381944.5
$ws.Range("H136").Value = 1755846.4
$ws.Range("I136").Value = 2382156.5
$ws.Range("K136").Value = 7146469.5
$ws.Range("M136").Value = -7143919.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1659.8
$ws.Range("I20").Value = 1483.6666
$ws.Range("J20").Value = 1924
$ws.Range("K20").Value = 1483.6666
$ws.Range("L20").Value = 1924
$ws.Range("M20").Value = -1236.6666
$ws.Range("N20").Value = -2418
$ws.Range("H64").Value = 2000.5
$ws.Range("J64").Value = 2000.5
$ws.Range("L64").Value = 2000.5
$ws.Range("N64").Value = -2450.5
$ws.Range("H67").Value = 2000.5
$ws.Range("J67").Value = 2000.5
$ws.Range("L67").Value = 2000.5
$ws.Range("N67").Value = -3560.5
$ws.Range("H86").Value = 3642.2856
$ws.Range("J86").Value = 5007
$ws.Range("L86").Value = 5007
$ws.Range("N86").Value = -7253
$ws.Range("H89").Value = 3642.2856
$ws.Range("J89").Value = 5007
$ws.Range("L89").Value = 25035
$ws.Range("N89").Value = -36267
$ws.Range("H134").Value = 1054058
$ws.Range("I134").Value = 1085173.8
$ws.Range("J134").Value = 917148.6
$ws.Range("K134").Value = 3255521.4
$ws.Range("L134").Value = 2751445.8
$ws.Range("M134").Value = -3252986.4
$ws.Range("N134").Value = -2756515.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3720
$ws.Range("I62").Value = 2400
$ws.Range("J62").Value = 4050
$ws.Range("K62").Value = 2400
$ws.Range("L62").Value = 4050
$ws.Range("M62").Value = -1776
$ws.Range("N62").Value = -5298
$ws.Range("H65").Value = 3720
$ws.Range("I65").Value = 2400
$ws.Range("J65").Value = 4050
$ws.Range("K65").Value = 12000
$ws.Range("L65").Value = 20250
$ws.Range("M65").Value = -8880
$ws.Range("N65").Value = -26490
$ws.Range("H132").Value = 32425810
$ws.Range("I132").Value = 37038964
$ws.Range("K132").Value = 111116892
$ws.Range("M132").Value = -111114362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 828.6
$ws.Range("J36").Value = 847.6667
$ws.Range("L36").Value = 2543.0001
$ws.Range("N36").Value = -2881.0001
$ws.Range("H129").Value = 1858.9445
$ws.Range("J129").Value = 3997.8572
$ws.Range("L129").Value = 11993.5716
$ws.Range("N129").Value = -21993.5716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 304102.88
$ws.Range("I80").Value = 467748.8
$ws.Range("J80").Value = 4085.3333
$ws.Range("K80").Value = 467748.8
$ws.Range("L80").Value = 4085.3333
$ws.Range("M80").Value = -466750.8
$ws.Range("N80").Value = -6081.3333
$ws.Range("H83").Value = 304102.88
$ws.Range("I83").Value = 467748.8
$ws.Range("J83").Value = 4085.3333
$ws.Range("K83").Value = 2338744
$ws.Range("L83").Value = 20426.6665
$ws.Range("M83").Value = -2333752
$ws.Range("N83").Value = -30410.6665
$ws.Range("H102").Value = 3045.838
$ws.Range("I102").Value = 2245.8462
$ws.Range("J102").Value = 4936.727
$ws.Range("K102").Value = 2245.8462
$ws.Range("L102").Value = 4936.727
$ws.Range("M102").Value = -623.8462
$ws.Range("N102").Value = -8180.727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 784.3333
$ws.Range("I16").Value = 506.58334
$ws.Range("J16").Value = 1339.8334
$ws.Range("K16").Value = 506.58334
$ws.Range("L16").Value = 1339.8334
$ws.Range("M16").Value = -336.58334
$ws.Range("N16").Value = -1679.8334
$ws.Range("H68").Value = 3309.818
$ws.Range("I68").Value = 3480.8
$ws.Range("K68").Value = 3480.8
$ws.Range("M68").Value = -2731.8
$ws.Range("H71").Value = 3309.818
$ws.Range("I71").Value = 3480.8
$ws.Range("K71").Value = 17404
$ws.Range("M71").Value = -13660
$ws.Range("H122").Value = 4849.591
$ws.Range("I122").Value = 4594.278
$ws.Range("J122").Value = 5998.5
$ws.Range("K122").Value = 13782.834
$ws.Range("L122").Value = 17995.5
$ws.Range("M122").Value = -11332.834
$ws.Range("N122").Value = -22895.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H132").Value = 3799893
$ws.Range("I132").Value = 4474408
$ws.Range("J132").Value = 5747.875
$ws.Range("K132").Value = 13423224
$ws.Range("L132").Value = 17243.625
$ws.Range("M132").Value = -13420694
$ws.Range("N132").Value = -22303.625
